# Update symbol list (crypto price/volume snapshot) per Thu Feb  2 19:23:48 UTC 2023 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'329.16"
$ws.Range("E2").Formula = "'6.25%"
$ws.Range("D3").Formula = "'40.67"
$ws.Range("E3").Formula = "'11.25%"
$ws.Range("D4").Formula = "'6.045"
$ws.Range("E4").Formula = "'18.26%"
$ws.Range("D5").Formula = "'0.08145"
$ws.Range("E5").Formula = "'5.50%"
$ws.Range("D6").Formula = "'4.605"
$ws.Range("E6").Formula = "'4.93%"
$ws.Range("D7").Formula = "'8.787"
$ws.Range("E7").Formula = "'5.72%"
$ws.Range("D8").Formula = "'1.973"
$ws.Range("E8").Formula = "'7.05%"
$ws.Range("E9").Formula = "'-0.01%"
$ws.Range("D10").Formula = "'0.9493"
$ws.Range("E10").Formula = "'2.83%"
$ws.Range("D11").Formula = "'0.1340"
$ws.Range("E11").Formula = "'16.80%"
$ws.Range("D12").Formula = "'0.1997"
$ws.Range("E12").Formula = "'6.67%"
$ws.Range("B13").Value = "MCDex"
$ws.Range("C13").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D13").Formula = "'9.616"
$ws.Range("E13").Formula = "'52.56%"
$ws.Range("B14").Value = "MandalaExchangeToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D14").Formula = "'0.09375"
$ws.Range("E14").Formula = "'6.88%"
$ws.Range("B15").Value = "BitrueCoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D15").Formula = "'0.03505"
$ws.Range("E15").Formula = "'3.97%"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").Formula = "'0.09613"
$ws.Range("E16").Formula = "'0.91%"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Formula = "'0.001315"
$ws.Range("E17").Formula = "'-4.27%"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Formula = "'0.006384"
$ws.Range("E18").Formula = "'9.55%"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Formula = "'3.353"
$ws.Range("E19").Formula = "'-0.09%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Formula = "'0.3545"
$ws.Range("E20").Formula = "'3.21%"
$ws.Range("D21").Formula = "'0.1423"
$ws.Range("E21").Formula = "'10.26%"
$ws.Range("E22").Formula = "'5.64%"
$ws.Range("D23").Formula = "'0.04437"
$ws.Range("E23").Formula = "'2.55%"
$ws.Range("D24").Formula = "'0.001262"
$ws.Range("E24").Formula = "'5.11%"
$ws.Range("D25").Formula = "'0.004382"
$ws.Range("E25").Formula = "'3.10%"
$ws.Range("D26").Formula = "'0.0001093"
$ws.Range("E26").Formula = "'-17.78%"
$ws.Range("E27").Formula = "'4.80%"
$ws.Range("D39").Formula = "'0.02491"
$ws.Range("E39").Formula = "'17.90%"
$ws.Range("D40").Formula = "'0.05304"
$ws.Range("E40").Formula = "'6.02%"
$ws.Range("D41").Formula = "'0.007474"
$ws.Range("E41").Formula = "'-0.28%"
$ws.Range("D42").Formula = "'0.1438"
$ws.Range("E42").Formula = "'6.78%"
$ws.Range("D43").Formula = "'0.009097"
$ws.Range("E43").Formula = "'8.35%"
$ws.Range("D44").Formula = "'0.002057"
$ws.Range("E44").Formula = "'-0.59%"
$ws.Range("D45").Formula = "'0.01056"
$ws.Range("E45").Formula = "'37.04%"
$ws.Range("D46").Formula = "'0.00006822"
$ws.Range("E46").Formula = "'8.10%"
$ws.Range("D47").Formula = "'0.00000000752"
$ws.Range("E47").Formula = "'0.33%"
$ws.Range("D48").Formula = "'0.003499"
$ws.Range("E48").Formula = "'22.12%"
$ws.Range("D49").Formula = "'0.001804"
$ws.Range("E49").Formula = "'6.81%"
$ws.Range("D50").Formula = "'0.00002106"
$ws.Range("E50").Formula = "'0.33%"
$ws.Range("D51").Formula = "'0.0002006"
$ws.Range("E51").Formula = "'0.33%"
